$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137
$ws.Range("A137").Value = 112101773
$ws.Range("B137").Value = 86357
$ws.Range("D137").Value = "NT"
$ws.Range("E137").Value = 4412
$ws.Range("F137").Value = "Äggvaxskivling"
$ws.Range("G137").Value = "Hygrophorus karstenii"
$ws.Range("H137").Value = "Sacc. & Cub."
$ws.Range("Q137").Value = 446984
$ws.Range("R137").Value = 7032942

# Row 138
$ws.Range("A138").Value = 112102104
$ws.Range("B138").Value = 90785
$ws.Range("D138").Value = "NT"
$ws.Range("E138").Value = 1968
$ws.Range("F138").Value = "Grantaggsvamp"
$ws.Range("G138").Value = "Bankera violascens"
$ws.Range("H138").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q138").Value = 446883
$ws.Range("R138").Value = 7032568

# Row 139
$ws.Range("A139").Value = 112111486
$ws.Range("B139").Value = 83072
$ws.Range("D139").Value = "NT"
$ws.Range("E139").Value = 5589
$ws.Range("F139").Value = "Rödbrun klubbdyna"
$ws.Range("G139").Value = "Trichoderma nybergianum"
$ws.Range("H139").Value = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
$ws.Range("Q139").Value = 446833
$ws.Range("R139").Value = 7032727

# Row 140
$ws.Range("A140").Value = 112103325
$ws.Range("B140").Value = 85434
$ws.Range("D140").Value = "NT"
$ws.Range("E140").Value = 3739
$ws.Range("F140").Value = "Persiljespindling"
$ws.Range("G140").Value = "Cortinarius sulfurinus"
$ws.Range("H140").Value = "Quél."
$ws.Range("Q140").Value = 446867
$ws.Range("R140").Value = 7032726

# Row 141
$ws.Range("A141").Value = 112102682
$ws.Range("B141").Value = 84941
$ws.Range("D141").Value = "VU"
$ws.Range("E141").Value = 275
$ws.Range("F141").Value = "Kejsarskivling"
$ws.Range("G141").Value = "Catathelasma imperiale"
$ws.Range("H141").Value = "(P.Karst.) Singer"
$ws.Range("Q141").Value = 447025
$ws.Range("R141").Value = 7032672

# Row 142
$ws.Range("A142").Value = 112104553
$ws.Range("B142").Value = 86357
$ws.Range("D142").Value = "NT"
$ws.Range("E142").Value = 4412
$ws.Range("F142").Value = "Äggvaxskivling"
$ws.Range("G142").Value = "Hygrophorus karstenii"
$ws.Range("H142").Value = "Sacc. & Cub."
$ws.Range("Q142").Value = 446688
$ws.Range("R142").Value = 7032560

# Row 143
$ws.Range("A143").Value = 112111498
$ws.Range("B143").Value = 88167
$ws.Range("D143").Value = "VU"
$ws.Range("E143").Value = 1599
$ws.Range("F143").Value = "Fjällfotad musseron"
$ws.Range("G143").Value = "Tricholoma olivaceotinctum"
$ws.Range("H143").Value = "Mort.Chr. & Heilm.-Claus."
$ws.Range("Q143").Value = 446860
$ws.Range("R143").Value = 7032743

# Row 144
$ws.Range("A144").Value = 112104266
$ws.Range("B144").Value = 89090
$ws.Range("D144").Value = "VU"
$ws.Range("E144").Value = 5747
$ws.Range("F144").Value = "Läderdoftande fingersvamp"
$ws.Range("G144").Value = "Ramaria safraniolens"
$ws.Range("H144").Value = "Christian"
$ws.Range("Q144").Value = 446732
$ws.Range("R144").Value = 7032598

# Row 145
$ws.Range("A145").Value = 112102200
$ws.Range("B145").Value = 89033
$ws.Range("D145").Value = "NT"
$ws.Range("E145").Value = 3286
$ws.Range("F145").Value = "Flattoppad klubbsvamp"
$ws.Range("G145").Value = "Clavariadelphus truncatus"
$ws.Range("H145").Value = "(Quél.) Donk"
$ws.Range("Q145").Value = 446961
$ws.Range("R145").Value = 7032566

# Row 146
$ws.Range("A146").Value = 112104573
$ws.Range("B146").Value = 88167
$ws.Range("D146").Value = "VU"
$ws.Range("E146").Value = 1599
$ws.Range("F146").Value = "Fjällfotad musseron"
$ws.Range("G146").Value = "Tricholoma olivaceotinctum"
$ws.Range("H146").Value = "Mort.Chr. & Heilm.-Claus."
$ws.Range("Q146").Value = 446696
$ws.Range("R146").Value = 7032530

# Row 147
$ws.Range("A147").Value = 112104547
$ws.Range("B147").Value = 89080
$ws.Range("D147").Value = "VU"
$ws.Range("E147").Value = 256335
$ws.Range("F147").Value = "Taggfingersvamp"
$ws.Range("G147").Value = "Ramaria karstenii"
$ws.Range("H147").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("Q147").Value = 446688
$ws.Range("R147").Value = 7032560

# Row 148
$ws.Range("A148").Value = 112101944
$ws.Range("B148").Value = 89033
$ws.Range("D148").Value = "NT"
$ws.Range("E148").Value = 3286
$ws.Range("F148").Value = "Flattoppad klubbsvamp"
$ws.Range("G148").Value = "Clavariadelphus truncatus"
$ws.Range("H148").Value = "(Quél.) Donk"
$ws.Range("Q148").Value = 446858
$ws.Range("R148").Value = 7032743

# Row 149
$ws.Range("A149").Value = 112111606
$ws.Range("B149").Value = 89033
$ws.Range("D149").Value = "NT"
$ws.Range("E149").Value = 3286
$ws.Range("F149").Value = "Flattoppad klubbsvamp"
$ws.Range("G149").Value = "Clavariadelphus truncatus"
$ws.Range("H149").Value = "(Quél.) Donk"
$ws.Range("Q149").Value = 446863
$ws.Range("R149").Value = 7032718

# Row 150
$ws.Range("A150").Value = 112104270
$ws.Range("B150").Value = 85434
$ws.Range("D150").Value = "NT"
$ws.Range("E150").Value = 3739
$ws.Range("F150").Value = "Persiljespindling"
$ws.Range("G150").Value = "Cortinarius sulfurinus"
$ws.Range("H150").Value = "Quél."
$ws.Range("Q150").Value = 446732
$ws.Range("R150").Value = 7032598

# Row 151
$ws.Range("A151").Value = 112102196
$ws.Range("B151").Value = 89090
$ws.Range("D151").Value = "VU"
$ws.Range("E151").Value = 5747
$ws.Range("F151").Value = "Läderdoftande fingersvamp"
$ws.Range("G151").Value = "Ramaria safraniolens"
$ws.Range("H151").Value = "Christian"
$ws.Range("Q151").Value = 446964
$ws.Range("R151").Value = 7032565
